$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Append the new "TOR" and "Steve" pathway rows (93-131) ---

$ws.Cells.Item(93, 1).Value = "TOR"
$ws.Cells.Item(93, 2).Value = "TOR1"
$ws.Cells.Item(93, 3).Value = "MGG_15156"
$ws.Cells.Item(94, 1).Value = "TOR"
$ws.Cells.Item(94, 2).Value = "LST8"
$ws.Cells.Item(94, 3).Value = "MGG_07284"
$ws.Cells.Item(95, 1).Value = "TOR"
$ws.Cells.Item(95, 2).Value = "MIP1"
$ws.Cells.Item(95, 3).Value = "MGG_02874"
$ws.Cells.Item(96, 1).Value = "TOR"
$ws.Cells.Item(96, 2).Value = "BIT61"
$ws.Cells.Item(96, 3).Value = "MGG_11443"
$ws.Cells.Item(97, 1).Value = "TOR"
$ws.Cells.Item(97, 2).Value = "MAPKAP1"
$ws.Cells.Item(97, 3).Value = "MGG_12955"
$ws.Cells.Item(98, 1).Value = "TOR"
$ws.Cells.Item(98, 2).Value = "spherulin-1b precursor"
$ws.Cells.Item(98, 3).Value = "MGG_15259"
$ws.Cells.Item(99, 1).Value = "TOR"
$ws.Cells.Item(99, 2).Value = "TSC11"
$ws.Cells.Item(99, 3).Value = "MGG_05284"
$ws.Cells.Item(100, 1).Value = "TOR"
$ws.Cells.Item(100, 2).Value = "IML1"
$ws.Cells.Item(100, 3).Value = "MGG_04160"
$ws.Cells.Item(101, 1).Value = "TOR"
$ws.Cells.Item(101, 2).Value = "SCH9"
$ws.Cells.Item(101, 3).Value = "MGG_14773"
$ws.Cells.Item(102, 1).Value = "TOR"
$ws.Cells.Item(102, 2).Value = "tor signaling pathway regulator"
$ws.Cells.Item(102, 3).Value = "MGG_01540"
$ws.Cells.Item(103, 1).Value = "TOR"
$ws.Cells.Item(103, 2).Value = "RPS6"
$ws.Cells.Item(103, 3).Value = "MGG_03236"
$ws.Cells.Item(104, 1).Value = "TOR"
$ws.Cells.Item(104, 2).Value = "YPK1"
$ws.Cells.Item(104, 3).Value = "MGG_06599"
$ws.Cells.Item(105, 1).Value = "TOR"
$ws.Cells.Item(105, 2).Value = "ORM1"
$ws.Cells.Item(105, 3).Value = "MGG_16259"
$ws.Cells.Item(106, 1).Value = "TOR"
$ws.Cells.Item(106, 2).Value = "LCB1"
$ws.Cells.Item(106, 3).Value = "MGG_00864"
$ws.Cells.Item(107, 1).Value = "TOR"
$ws.Cells.Item(107, 2).Value = "LCB2"
$ws.Cells.Item(107, 3).Value = "MGG_05197"
$ws.Cells.Item(108, 1).Value = "TOR"
$ws.Cells.Item(108, 2).Value = "phosphoinositide phosphatase"
$ws.Cells.Item(108, 3).Value = "MGG_11279"
$ws.Cells.Item(109, 1).Value = "TOR"
$ws.Cells.Item(109, 2).Value = "PP2a"
$ws.Cells.Item(109, 3).Value = "MGG_05637"
$ws.Cells.Item(110, 1).Value = "TOR"
$ws.Cells.Item(110, 2).Value = "LAC1"
$ws.Cells.Item(110, 3).Value = "MGG_05189"
$ws.Cells.Item(111, 1).Value = "TOR"
$ws.Cells.Item(111, 2).Value = "LAG1"
$ws.Cells.Item(111, 3).Value = "MGG_03090"
$ws.Cells.Item(112, 1).Value = "TOR"
$ws.Cells.Item(112, 2).Value = "ELO1"
$ws.Cells.Item(112, 3).Value = "MGG_07280"
$ws.Cells.Item(113, 1).Value = "TOR"
$ws.Cells.Item(113, 2).Value = "FPK1"
$ws.Cells.Item(113, 3).Value = "MGG_07012"
$ws.Cells.Item(114, 1).Value = "TOR"
$ws.Cells.Item(114, 2).Value = "MgAPT2"
$ws.Cells.Item(114, 3).Value = "MGG_02767"
$ws.Cells.Item(115, 1).Value = "TOR"
$ws.Cells.Item(115, 2).Value = "GIN4"
$ws.Cells.Item(115, 3).Value = "MGG_02810"
$ws.Cells.Item(116, 1).Value = "TOR"
$ws.Cells.Item(116, 2).Value = "GPD1"
$ws.Cells.Item(116, 3).Value = "MGG_00067"
$ws.Cells.Item(117, 1).Value = "TOR"
$ws.Cells.Item(117, 2).Value = "mip family channel protein"
$ws.Cells.Item(117, 3).Value = "MGG_13615"
$ws.Cells.Item(118, 1).Value = "TOR"
$ws.Cells.Item(118, 2).Value = "ph domain-containing protein"
$ws.Cells.Item(118, 3).Value = "MGG_03558"
$ws.Cells.Item(119, 1).Value = "TOR"
$ws.Cells.Item(119, 2).Value = "ph domain-containing protein"
$ws.Cells.Item(119, 3).Value = "MGG_01047"
$ws.Cells.Item(120, 1).Value = "TOR"
$ws.Cells.Item(120, 2).Value = "PKC1"
$ws.Cells.Item(120, 3).Value = "MGG_08689"
$ws.Cells.Item(121, 1).Value = "TOR"
$ws.Cells.Item(121, 2).Value = "RHO1"
$ws.Cells.Item(121, 3).Value = "MGG_07176"
$ws.Cells.Item(122, 1).Value = "TOR"
$ws.Cells.Item(122, 2).Value = "rho1 guanine nucleotide exchange factor 1"
$ws.Cells.Item(122, 3).Value = "MGG_03064"
$ws.Cells.Item(123, 1).Value = "TOR"
$ws.Cells.Item(123, 2).Value = "MPS1"
$ws.Cells.Item(123, 3).Value = "MGG_04943"
$ws.Cells.Item(124, 1).Value = "TOR"
$ws.Cells.Item(124, 2).Value = "KSG1"
$ws.Cells.Item(124, 3).Value = "MGG_01795"
$ws.Cells.Item(125, 1).Value = "TOR"
$ws.Cells.Item(125, 2).Value = "CPKA"
$ws.Cells.Item(125, 3).Value = "MGG_06368"
$ws.Cells.Item(126, 1).Value = "TOR"
$ws.Cells.Item(126, 2).Value = "EHS1"
$ws.Cells.Item(126, 3).Value = "MGG_12128"
$ws.Cells.Item(127, 1).Value = "TOR"
$ws.Cells.Item(127, 2).Value = "CMKK2/MoTos3"
$ws.Cells.Item(127, 3).Value = "MGG_06421"
$ws.Cells.Item(128, 1).Value = "Steve"
$ws.Cells.Item(128, 2).Value = "wdA (AN1056)"
$ws.Cells.Item(128, 3).Value = "MGG_06968"
$ws.Cells.Item(129, 1).Value = "Steve"
$ws.Cells.Item(129, 2).Value = "mcnC/Def1 (AN2871)"
$ws.Cells.Item(129, 3).Value = "MGG_00124"
$ws.Cells.Item(130, 1).Value = "Steve"
$ws.Cells.Item(130, 2).Value = "Tubulin binding Cofactor A (TBCA)  (AN6176)"
$ws.Cells.Item(130, 3).Value = "MGG_09890"
$ws.Cells.Item(131, 1).Value = "Steve"
$ws.Cells.Item(131, 2).Value = " Elongation Factor 2 (AN6330)"
$ws.Cells.Item(131, 3).Value = "MGG_01742"

# Column A on these rows uses the same plain-Arial style as the rest of the table
# (row 92 is a representative existing data row); copy formats only so the new
# shared-string values we just wrote are not touched.
$ws.Range("A92").Copy()
$ws.Range("A93:A131").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Match the saved selection/scroll state from the edited workbook
$ws.Range("A131").Select()

Write-Output "done"